$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: per-person cost 35 -> 30, with a note about the change ---
$ws.Range("B10").Value = 30
$ws.Range("D10").Value = '(changed from $35)'

# --- Row 33: groceries cost formula gained an extra $6.45 item ---
$ws.Range("D33").Formula = '=66.15+6.45-5-2.5-2.99'

# --- Row 34: gas line re-attributed from Schuster to Pauly, and the
#     formula total replaced with the actual (static) receipt amount ---
$ws.Range("B34").Value = "Pauly"
$ws.Range("D34").Value = 116.35

# --- Row 50 (Stephan SA): cost now references $B$10 instead of the
#     literal 35; D50 becomes a live ABS() formula; E50 gets a note ---
$ws.Range("B50").Formula = '=3*$B$10-D33-39'
$ws.Range("D50").Formula = '=ABS(B50)'
$ws.Range("E50").Value = "credit to Stephan account"

# --- Row 51 (Schuster SA): formula no longer nets out D34 (gas is now
#     Pauly's line); D51 becomes the static "debit" label; E51 cleared ---
$ws.Range("B51").Formula = '=2*$B$10-39'
$ws.Range("D51").Value = "debit"
$ws.Range("E51").Value = ""

# --- Rows 52-56: literal 35 dues become live references to $B$10 ---
$ws.Range("B52").Formula = '=$B$10'
$ws.Range("B53").Formula = '=$B$10'
$ws.Range("B54").Formula = '=$B$10'
$ws.Range("B55").Formula = '=$B$10'
$ws.Range("B56").Formula = '=$B$10'

# stray formatting touches that appear in the new workbook
$ws.Range("I52").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("G54").NumberFormat = "$#,##0_);[Red]($#,##0)"

# --- Row 57: now holds the Pauly gas credit/debit entries ---
$ws.Range("A57").Value = "Pauly SA"
$ws.Range("B57").Formula = '=-D34-2*5'
$ws.Range("D57").Formula = '=ABS(B57)'
$ws.Range("E57").Value = "credit to Pauly account"

# --- Row 59 (new): actual total line ---
$ws.Range("D59").Value = "actual total"
$ws.Range("E59").Formula = '=SUM(B49:B57)+B38'
$ws.Range("E59").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

$ws.Range("A59").Select()
